# Update the localization-status report:
#   - "Ready for handoff" -> "In Translation" wherever it appears
#     (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - Narrow the now-shorter "Status" columns to match (Overview cols E/F,
#     zh-cn/de-de col C) as close as the host's column-width granularity
#     allows.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Closest width (in Excel "characters" ColumnWidth units) this host can
# store that lands nearest the narrower target column width used by the
# original report generator.
$newColWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
                $ws.Columns.Item($c).ColumnWidth = $newColWidth
            }
        }
    }
}
